# edit.ps1 — reproduces the "Add files via upload" revision of Recursos.xlsx
# (new Reservas bookings, Horas list reflowed, Mantenimientos layout change).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Reservas" sheet — replace the old booking rows (2-19) with the six new
#    bookings from the refreshed export. Keep the header row (1) untouched.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Reservas")

# Drop the now-obsolete tail rows (8-19) first so the remaining six data
# rows (2-7) line up with the new dimension (A1:G7).
$ws.Range("A8:G19").EntireRow.Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = 45852
$ws.Range("B2").Value = 0.3333333333333333
$ws.Range("C2").Value = 0.3958333333333333
$ws.Range("D2").Value = "ALEXIS FERNANDO SARZOSA VÉLIZ"
$ws.Range("E2").Value = "3° BÁSICO A"
$ws.Range("F2").Value = "ENLACE BASICA"
$ws.Range("G2").Value = "umaximo"

# Row 3
$ws.Range("A3").Value = 45853
$ws.Range("B3").Value = 0.3333333333333333
$ws.Range("C3").Value = 0.3958333333333333
$ws.Range("D3").Value = "ANDREA VICENTA DOSSI SAAVEDRA"
$ws.Range("E3").Value = " 2° BÁSICO B"
$ws.Range("F3").Value = "MICROFONOS"
$ws.Range("G3").Value = "lentes"

# Row 4
$ws.Range("A4").Value = 45853
$ws.Range("B4").Value = 0.3333333333333333
$ws.Range("C4").Value = 0.3958333333333333
$ws.Range("D4").Value = "ANDREA VICENTA DOSSI SAAVEDRA"
$ws.Range("E4").Value = " 2° BÁSICO B"
$ws.Range("F4").Value = "MOVIL 3"
$ws.Range("G4").Value = "vr"

# Row 5
$ws.Range("A5").Value = 45853
$ws.Range("B5").Value = 0.3333333333333333
$ws.Range("C5").Value = 0.3958333333333333
$ws.Range("D5").Value = "ANDREA VICENTA DOSSI SAAVEDRA"
$ws.Range("E5").Value = " 3° BÁSICO A"
$ws.Range("F5").Value = "LAPICES 3D"
$ws.Range("G5").Value = "test 2334"

# Row 6
$ws.Range("A6").Value = 45852
$ws.Range("B6").Value = 0.3333333333333333
$ws.Range("C6").Value = 0.3958333333333333
$ws.Range("D6").Value = "ALEXIS FERNANDO SARZOSA VÉLIZ"
$ws.Range("E6").Value = " 2° BÁSICO A"
$ws.Range("F6").Value = "LAPICES 3D"
$ws.Range("G6").Value = "alex"

# Row 7
$ws.Range("A7").Value = 45853
$ws.Range("B7").Value = 0.3333333333333333
$ws.Range("C7").Value = 0.3958333333333333
$ws.Range("D7").Value = "BRENDA LISSETTE PÉREZ ESCOBAR"
$ws.Range("E7").Value = " 2° MEDIO A"
$ws.Range("F7").Value = "LENTES VR"
$ws.Range("G7").Value = "alex2"

$ws.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) "Horas" sheet — the list used to start on row 2 (row 1 was a stray
#    blank); remove that leading blank row so the schedule labels shift up
#    to A1:A14.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Horas")
$ws.Range("A1").EntireRow.Delete() | Out-Null
$ws.Range("A1:A14").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) "Recursos" sheet (the resource lookup list) — no data change, just the
#    cursor was left on D6 when the file was last saved.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Recursos")
$ws.Range("D6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) "Mantenimientos" sheet — swap the FechaFin/HoraInicio columns (HoraInicio
#    now comes right after FechaInicio) and refresh the sample maintenance
#    row to the new dates.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mantenimientos")

# Header row: C was "FechaFin", D was "HoraInicio" — swap their labels so the
# columns read Recurso | FechaInicio | HoraInicio | FechaFin | HoraFin.
$ws.Range("C1").Value = "HoraInicio"
$ws.Range("D1").Value = "FechaFin"

# Match the bold/boxed header look used elsewhere in the workbook.
$hdr = $ws.Range("A1:E1")
$hdr.Font.Bold = $True
$hdr.Font.Name = "Calibri"
$hdr.Font.Size = 11
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data row: FechaInicio/FechaFin both move to 2025-07-15 (45853); HoraInicio
# becomes plain text in its new column C, FechaFin keeps the date format in
# its new column D.
$ws.Range("B2").Value = 45853
$ws.Range("C2").ClearFormats() | Out-Null
$ws.Range("C2").Value = "08:00:00"
$ws.Range("D2").Value = 45853
$ws.Range("D2").NumberFormat = $ws.Range("B2").NumberFormat

# New column widths for FechaInicio (B) and FechaFin (D).
$ws.Columns.Item(2).ColumnWidth = 11
$ws.Columns.Item(4).ColumnWidth = 13.5

$ws.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5) Leave the workbook on the "Reservas" tab, matching the refreshed file.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Reservas").Activate() | Out-Null
